$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 334.1
$ws.Range("I38").Value = 260.1111
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 780.3333
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = -408.3333
$ws.Range("N38").Value = -3744

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5852.353
$ws.Range("J40").Value = 6785
$ws.Range("L40").Value = 6785
$ws.Range("N40").Value = -7135

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2643.7144
$ws.Range("I100").Value = 2900
$ws.Range("J100").Value = 2003
$ws.Range("K100").Value = 2900
$ws.Range("L100").Value = 2003
$ws.Range("M100").Value = -2359
$ws.Range("N100").Value = -3085

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1645.5555
$ws.Range("I137").Value = 1261.7142
$ws.Range("J137").Value = 2989
$ws.Range("K137").Value = 3785.1426
$ws.Range("L137").Value = 8967
$ws.Range("M137").Value = -1235.1426
$ws.Range("N137").Value = -14067

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2063.7334
$ws.Range("I138").Value = 1022.2857
$ws.Range("J138").Value = 2975
$ws.Range("K138").Value = 3066.8571
$ws.Range("L138").Value = 8925
$ws.Range("M138").Value = 2073.1429
$ws.Range("N138").Value = -19205

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 52000
$ws.Range("J29").Value = 52000
$ws.Range("L29").Value = 52000
$ws.Range("N29").Value = -52616

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2004.7142
$ws.Range("I74").Value = 2004.7142
$ws.Range("K74").Value = 2004.7142
$ws.Range("M74").Value = -1130.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2004.7142
$ws.Range("I77").Value = 2004.7142
$ws.Range("K77").Value = 10023.571
$ws.Range("M77").Value = -5655.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2257
$ws.Range("I102").Value = 2257
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2257
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -635
$ws.Range("N102").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 50000
$ws.Range("J104").Value = 50000
$ws.Range("L104").Value = 50000
$ws.Range("N104").Value = -56988

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2336.0908
$ws.Range("I110").Value = 1242.5714
$ws.Range("K110").Value = 1242.5714
$ws.Range("M110").Value = 802.4286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 347.66666
$ws.Range("I64").Value = 347.5
$ws.Range("J64").Value = 348
$ws.Range("K64").Value = 347.5
$ws.Range("L64").Value = 348
$ws.Range("M64").Value = -122.5
$ws.Range("N64").Value = -798

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 347.66666
$ws.Range("I67").Value = 347.5
$ws.Range("J67").Value = 348
$ws.Range("K67").Value = 347.5
$ws.Range("L67").Value = 348
$ws.Range("M67").Value = 432.5
$ws.Range("N67").Value = -1908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 10000
$ws.Range("J88").Value = 10000
$ws.Range("L88").Value = 10000
$ws.Range("N88").Value = -10812

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H91").Value = 10000
$ws.Range("J91").Value = 10000
$ws.Range("L91").Value = 10000
$ws.Range("N91").Value = -12808

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1121.4546
$ws.Range("I107").Value = 1121.4546
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1121.4546
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 798.5454
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 191.16667
$ws.Range("I5").Value = 189.8
$ws.Range("J5").Value = 198
$ws.Range("K5").Value = 189.8
$ws.Range("L5").Value = 198
$ws.Range("M5").Value = -77.80000000000001
$ws.Range("N5").Value = -422

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2521.0908
$ws.Range("I31").Value = 2610.5715
$ws.Range("J31").Value = 2364.5
$ws.Range("K31").Value = 2610.5715
$ws.Range("L31").Value = 2364.5
$ws.Range("M31").Value = -2315.5715
$ws.Range("N31").Value = -2954.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 4357.143
$ws.Range("J33").Value = 6000
$ws.Range("L33").Value = 6000
$ws.Range("N33").Value = -6758

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2521.0908
$ws.Range("I34").Value = 2610.5715
$ws.Range("J34").Value = 2364.5
$ws.Range("K34").Value = 2610.5715
$ws.Range("L34").Value = 2364.5
$ws.Range("M34").Value = -2408.5715
$ws.Range("N34").Value = -2768.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2696.0715
$ws.Range("I58").Value = 1632.8334
$ws.Range("J58").Value = 3493.5
$ws.Range("K58").Value = 1632.8334
$ws.Range("L58").Value = 3493.5
$ws.Range("M58").Value = -1429.8334
$ws.Range("N58").Value = -3899.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1965.2
$ws.Range("I122").Value = 1978
$ws.Range("K122").Value = 5934
$ws.Range("M122").Value = -3484

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2696.0715
$ws.Range("I136").Value = 1632.8334
$ws.Range("J136").Value = 3493.5
$ws.Range("K136").Value = 4898.5002
$ws.Range("L136").Value = 10480.5
$ws.Range("M136").Value = -2348.5002
$ws.Range("N136").Value = -15580.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 4850
$ws.Range("J132").Value = 4850
$ws.Range("L132").Value = 43650
$ws.Range("N132").Value = -48710

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1728.5
$ws.Range("I102").Value = 1804.8334
$ws.Range("J102").Value = 1499.5
$ws.Range("K102").Value = 1804.8334
$ws.Range("L102").Value = 1499.5
$ws.Range("M102").Value = -182.8334
$ws.Range("N102").Value = -4743.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 11588.714
$ws.Range("J126").Value = 13249.5
$ws.Range("L126").Value = 39748.5
$ws.Range("N126").Value = -44688.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4235.2856
$ws.Range("I46").Value = 3549
$ws.Range("J46").Value = 4750
$ws.Range("K46").Value = 3549
$ws.Range("L46").Value = 4750
$ws.Range("M46").Value = -3361
$ws.Range("N46").Value = -5126

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 686.125
$ws.Range("I81").Value = 712.7143
$ws.Range("J81").Value = 500
$ws.Range("K81").Value = 1425.4286
$ws.Range("L81").Value = 1000
$ws.Range("M81").Value = -364.4286
$ws.Range("N81").Value = -3122

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 686.125
$ws.Range("I84").Value = 712.7143
$ws.Range("J84").Value = 500
$ws.Range("K84").Value = 7127.143
$ws.Range("L84").Value = 5000
$ws.Range("M84").Value = -1823.143
$ws.Range("N84").Value = -15608

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 23622.875
$ws.Range("J126").Value = 36999.6
$ws.Range("L126").Value = 110998.8
$ws.Range("N126").Value = -115938.8
